$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume values
$updates = @(
    @{ Addr = "D2"; Value = "28.636.03" },
    @{ Addr = "E2"; Value = "  +2.27%  " },
    @{ Addr = "D3"; Value = "1.867.52" },
    @{ Addr = "E3"; Value = "  +2.14%  " },
    @{ Addr = "E4"; Value = "  +0.15%  " },
    @{ Addr = "D5"; Value = "323.90" },
    @{ Addr = "E5"; Value = "  -0.23%  " },
    @{ Addr = "D6"; Value = "1.001" },
    @{ Addr = "E6"; Value = "  +0.00%  " },
    @{ Addr = "D7"; Value = "0.4609" },
    @{ Addr = "E7"; Value = "  -0.77%  " },
    @{ Addr = "D8"; Value = "0.3876" },
    @{ Addr = "E8"; Value = "  +0.44%  " },
    @{ Addr = "D9"; Value = "0.07869" },
    @{ Addr = "E9"; Value = "  +0.19%  " },
    @{ Addr = "D10"; Value = "0.9751" },
    @{ Addr = "E10"; Value = "  +1.86%  " },
    @{ Addr = "D11"; Value = "21.96" },
    @{ Addr = "E11"; Value = "  +0.60%  " },
    @{ Addr = "D12"; Value = "1.871.75" },
    @{ Addr = "E12"; Value = "  +1.65%  " },
    @{ Addr = "D13"; Value = "6.996" },
    @{ Addr = "E13"; Value = "  +1.58%  " },
    @{ Addr = "D14"; Value = "5.695" },
    @{ Addr = "E14"; Value = "  +0.47%  " },
    @{ Addr = "D15"; Value = "0.06949" },
    @{ Addr = "E15"; Value = "  +1.21%  " },
    @{ Addr = "D16"; Value = "88.26" },
    @{ Addr = "E16"; Value = "  +1.30%  " },
    @{ Addr = "D17"; Value = "1.001" },
    @{ Addr = "E17"; Value = "  +0.07%  " },
    @{ Addr = "D18"; Value = "0.00001001" },
    @{ Addr = "E18"; Value = "  +0.98%  " },
    @{ Addr = "D19"; Value = "16.82" },
    @{ Addr = "E19"; Value = "  +1.55%  " },
    @{ Addr = "E20"; Value = "  +0.15%  " },
    @{ Addr = "D21"; Value = "28.618.23" },
    @{ Addr = "E21"; Value = "  +2.11%  " },
    @{ Addr = "D22"; Value = "5.267" },
    @{ Addr = "E22"; Value = "  -0.81%  " },
    @{ Addr = "D23"; Value = "11.08" },
    @{ Addr = "E23"; Value = "  +1.06%  " },
    @{ Addr = "D24"; Value = "2.103" },
    @{ Addr = "E24"; Value = "  +0.56%  " },
    @{ Addr = "D25"; Value = "2.080.71" },
    @{ Addr = "E25"; Value = "  +1.09%  " },
    @{ Addr = "D26"; Value = "152.48" },
    @{ Addr = "E26"; Value = "  -0.64%  " },
    @{ Addr = "D27"; Value = "19.26" },
    @{ Addr = "E27"; Value = "  +0.86%  " },
    @{ Addr = "D28"; Value = "5.868" },
    @{ Addr = "E28"; Value = "  +3.76%  " },
    @{ Addr = "D29"; Value = "1.984" },
    @{ Addr = "E29"; Value = "  +1.58%  " },
    @{ Addr = "D30"; Value = "119.26" },
    @{ Addr = "E30"; Value = "  +1.54%  " },
    @{ Addr = "D31"; Value = "0.09324" },
    @{ Addr = "E31"; Value = "  +1.13%  " },
    @{ Addr = "D32"; Value = "0.9191" },
    @{ Addr = "E32"; Value = "  -1.28%  " },
    @{ Addr = "D33"; Value = "5.271" },
    @{ Addr = "E33"; Value = "  +0.30%  " },
    @{ Addr = "D34"; Value = "1.333" },
    @{ Addr = "E34"; Value = "  +1.26%  " },
    @{ Addr = "D35"; Value = "3.316" },
    @{ Addr = "E35"; Value = "  +0.75%  " },
    @{ Addr = "D36"; Value = "0.05791" },
    @{ Addr = "E36"; Value = "  -0.66%  " },
    @{ Addr = "D37"; Value = "1.150" },
    @{ Addr = "E37"; Value = "  +0.97%  " },
    @{ Addr = "D38"; Value = "0.02071" },
    @{ Addr = "E38"; Value = "  -1.99%  " },
    @{ Addr = "D39"; Value = "7.685" },
    @{ Addr = "E39"; Value = "  -1.34%  " },
    @{ Addr = "D40"; Value = "0.5611" },
    @{ Addr = "E40"; Value = "  +0.60%  " },
    @{ Addr = "D41"; Value = "0.1782" },
    @{ Addr = "E41"; Value = "  +1.34%  " },
    @{ Addr = "D42"; Value = "9.773" },
    @{ Addr = "E42"; Value = "  -0.73%  " },
    @{ Addr = "D43"; Value = "0.07216" },
    @{ Addr = "D44"; Value = "11.64" },
    @{ Addr = "E44"; Value = "  +0.49%  " },
    @{ Addr = "D45"; Value = "0.5286" },
    @{ Addr = "E45"; Value = "  +0.54%  " },
    @{ Addr = "D46"; Value = "2.149" },
    @{ Addr = "E46"; Value = "  +1.58%  " },
    @{ Addr = "E47"; Value = "  +1.57%  " },
    @{ Addr = "D48"; Value = "1.836" },
    @{ Addr = "E48"; Value = "  +0.62%  " },
    @{ Addr = "D49"; Value = "112.83" },
    @{ Addr = "E49"; Value = "  +0.29%  " },
    @{ Addr = "E50"; Value = "  +3.44%  " },
    @{ Addr = "D51"; Value = "1.001" },
    @{ Addr = "E51"; Value = "  +0.09%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
